$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new file/table mapping row for the FY2022 Q1 revised pw_worksites form.
$ws.Range("A9").Value = "PW_Worksites_FY2022_Q1_revised_form_h2b.xlsx"
$ws.Range("B9").Value = "pw_worksites_new"

# Add a sample-marker row a few rows further down (row 12, column A only).
$ws.Range("A12").Value = "THIS IS A SAMPLE"

# Reflect the post-edit UI selection state (rows 2-7 selected across the full row width).
$ws.Range("A2:XFD7").Select()
